$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 ---
$ws.Range("A9").Value = 112395204
$ws.Range("B9").Value = 89317
$ws.Range("C9").Value = "Ovaliderad"
$ws.Range("D9").Value = "LC"
$ws.Range("E9").Value = 3215
$ws.Range("F9").Value = "Rödgul trumpetsvamp"
$ws.Range("G9").Value = "Craterellus lutescens"
$ws.Range("H9").Value = "(Fr.) Fr."
$ws.Range("P9").Value = "Sjöändan, Vrm"
$ws.Range("Q9").Value = 425634
$ws.Range("R9").Value = 6614497
$ws.Range("S9").Value = 10
$ws.Range("T9").Value = "Värmland"
$ws.Range("U9").Value = "Karlstad"
$ws.Range("V9").Value = "Värmland"
$ws.Range("W9").Value = "Nyed"

$ws.Range("Y9").NumberFormat = "@"
$ws.Range("Y9").Value = "2023-08-24"
$ws.Range("AA9").NumberFormat = "@"
$ws.Range("AA9").Value = "2023-08-24"

$ws.Range("AD9").Value = $false
$ws.Range("AE9").Value = $false
$ws.Range("AG9").Value = $false

$ws.Range("AW9").Value = "Jan Rees"
$ws.Range("AX9").Value = "Jan Rees"

# --- Row 10 ---
$ws.Range("A10").Value = 112395203
$ws.Range("B10").Value = 90821
$ws.Range("C10").Value = "Ovaliderad"
$ws.Range("D10").Value = "LC"
$ws.Range("E10").Value = 5964
$ws.Range("F10").Value = "Fjällig taggsvamp s.str."
$ws.Range("G10").Value = "Sarcodon imbricatus s.str."
$ws.Range("H10").Value = "(L.:Fr.) P.Karst."
$ws.Range("P10").Value = "Sjöändan, Vrm"
$ws.Range("Q10").Value = 425883
$ws.Range("R10").Value = 6614337
$ws.Range("S10").Value = 10
$ws.Range("T10").Value = "Värmland"
$ws.Range("U10").Value = "Karlstad"
$ws.Range("V10").Value = "Värmland"
$ws.Range("W10").Value = "Nyed"

$ws.Range("Y10").NumberFormat = "@"
$ws.Range("Y10").Value = "2023-08-24"
$ws.Range("AA10").NumberFormat = "@"
$ws.Range("AA10").Value = "2023-08-24"

$ws.Range("AD10").Value = $false
$ws.Range("AE10").Value = $false
$ws.Range("AG10").Value = $false

$ws.Range("AW10").Value = "Jan Rees"
$ws.Range("AX10").Value = "Jan Rees"
